$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Use a scratch cell far outside the used range to produce the new date
# string as a plain text formula result ("2023-11-06"), then copy/paste
# its *value* (xlPasteValues) into place. Pasting values only transfers
# the stored content - not the cell's number format - so the destination
# keeps its original (default) style instead of Excel re-sniffing the
# "2023-11-06" text and auto-converting it to a date serial number.
$ws.Range("Z1").Formula = "=""2023-11-06"""
$ws.Range("Z1").Copy()
$ws.Range("B2").PasteSpecial(-4163)
$ws.Range("B3").PasteSpecial(-4163)
$ws.Range("Z1").ClearContents()

# asistanceType column: plain text, no date-sniffing concern.
$ws.Range("D2").Value = "FALTA"
$ws.Range("D3").Value = "FALTA"
